# Applies the commit "Results in .docx and .txt for Special 12 panel regs":
# appends two new SourceCode paragraphs at the end of the document, reproducing
# the R "panel_est(Form_special_12_all, data_special_12_all)" call and its
# plm() regression output, styled the same way as the existing R chunks.

$d = $word.ActiveDocument

# wdLineBreak -- manual line break (renders as <w:br/> rather than a new <w:p/>)
$wdLineBreak = 6

# Data describing the two new paragraphs to append, each a list of "runs":
#   @{ Type = "text"; Style = <character style name>; Text = <literal text> }
#   @{ Type = "br" }                                   (manual line break)
$newParagraphs = @(
    @{
        PStyle = 'SourceCode'
        Runs = @(
            @{ Type = "text"; Style = 'CommentTok'; Text = '# Special 12 all assets together' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'KeywordTok'; Text = 'panel_est' },
            @{ Type = "text"; Style = 'NormalTok'; Text = '(Form_special_12_all, data_special_12_all)' }
        )
    },
    @{
        PStyle = 'SourceCode'
        Runs = @(
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Oneway (individual) effect Within Model' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Call:' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'plm::plm(formula = form, data = data_matrix, model = mdl, index = ind)' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Balanced Panel: n = 12, T = 21, N = 252' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Residuals:' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = '     Min.   1st Qu.    Median   3rd Qu.      Max. ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = '-28.07825  -5.71872  -0.63342   5.04082  36.14111 ' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Coefficients:' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = '           Estimate Std. Error t value  Pr(>|t|)    ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'TED       -5.932488   3.611108 -1.6428   0.10176    ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'VIX        0.146661   0.081247  1.8051   0.07235 .  ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'SENT      -0.727099   0.512025 -1.4200   0.15693    ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'FEDFUNDS   0.471201   0.263635  1.7873   0.07518 .  ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'INTERNET  -1.612662   0.183372 -8.7945 3.248e-16 ***' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'ERM      -20.308476   3.313825 -6.1284 3.754e-09 ***' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Euro      -1.548770   2.987810 -0.5184   0.60470    ' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = '---' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Signif. codes:  0 ''***'' 0.001 ''**'' 0.01 ''*'' 0.05 ''.'' 0.1 '' '' 1' },
            @{ Type = "br" },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Total Sum of Squares:    76282' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Residual Sum of Squares: 17456' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'R-Squared:      0.77117' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'Adj. R-Squared: 0.75349' },
            @{ Type = "br" },
            @{ Type = "text"; Style = 'VerbatimChar'; Text = 'F-statistic: 112.174 on 7 and 233 DF, p-value: < 2.22e-16' }
        )
    }
)

foreach ($paraSpec in $newParagraphs) {
    # Append a brand-new paragraph after the current last paragraph, carrying
    # the "SourceCode" paragraph style used by the rest of the R chunks.
    $d.Paragraphs.Last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Style = $paraSpec.PStyle

    foreach ($run in $paraSpec.Runs) {
        # Re-fetch the paragraph's range fresh on every run so the insertion
        # point reflects what was just inserted (a cached/reused Range object
        # does not reliably reposition itself across InsertBreak calls).
        $cur = $newPara.Range
        $cur.Collapse(0)
        if ($run.Type -eq "br") {
            $cur.InsertBreak($wdLineBreak)
        } else {
            $cur.InsertAfter($run.Text)
            $cur.Style = $run.Style
        }
    }
}

Write-Output "Appended $($newParagraphs.Count) paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
